$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the room name in A3: "Salle Etienne d'Orves" -> "Salle d'Estienne d'Orves"
$ws.Range("A3").Value = "Salle d'Estienne d'Orves"

# Widen column A so the longer room name keeps fitting (mirrors Excel's own
# best-fit recalculation after the text changed).
$ws.Columns("A").ColumnWidth = 19.33

# Reflect the active selection left after the edit
$ws.Range("A3").Select() | Out-Null
